# Apply the updated cryptocurrency price/volume snapshot described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'95.206.22"
$ws.Range("E2").Value = "  -1.34%  "

$ws.Range("D3").Value = "'3.575.37"
$ws.Range("E3").Value = "  -0.48%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'235.41"
$ws.Range("E5").Value = "  -1.57%  "

$ws.Range("D6").Value = "'658.59"
$ws.Range("E6").Value = "  +3.17%  "

$ws.Range("D7").Value = "'1.48"
$ws.Range("E7").Value = "  -0.26%  "

$ws.Range("D8").Value = "'0.399"
$ws.Range("E8").Value = "  -1.04%  "

$ws.Range("E9").Value = "  +0.06%  "

$ws.Range("D10").Value = "'1.00"
$ws.Range("E10").Value = "  -1.49%  "

$ws.Range("D11").Value = "'3.572.92"
$ws.Range("E11").Value = "  -0.49%  "

$ws.Range("E12").Value = "  +1.15%  "

$ws.Range("D13").Value = "'42.44"
$ws.Range("E13").Value = "  -1.65%  "

$ws.Range("E14").Value = "  +1.05%  "

$ws.Range("D15").Value = "'4.240.78"
$ws.Range("E15").Value = "  -0.61%  "

$ws.Range("D16").Value = "'95.039.72"
$ws.Range("E16").Value = "  -1.32%  "

$ws.Range("E17").Value = "  -0.38%  "

$ws.Range("D18").Value = "'3.580.10"
$ws.Range("E18").Value = "  -0.65%  "

$ws.Range("D19").Value = "'7.74"
$ws.Range("E19").Value = "  -3.41%  "

$ws.Range("D20").Value = "'12.62"
$ws.Range("E20").Value = "  -5.04%  "

$ws.Range("D21").Value = "'17.85"
$ws.Range("E21").Value = "  -1.24%  "

$ws.Range("D22").Value = "'3.46"
$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("D23").Value = "'508.87"
$ws.Range("E23").Value = "  -1.56%  "

$ws.Range("D24").Value = "'0.479"
$ws.Range("E24").Value = "  -3.57%  "

$ws.Range("D25").Value = "'6.86"
$ws.Range("E25").Value = "  +2.44%  "

$ws.Range("E26").Value = "  -1.69%  "

$ws.Range("D27").Value = "'95.20"
$ws.Range("E27").Value = "  -2.08%  "

$ws.Range("D28").Value = "'12.70"
$ws.Range("E28").Value = "  +1.95%  "

$ws.Range("D29").Value = "'3.766.09"
$ws.Range("E29").Value = "  -0.69%  "

$ws.Range("D30").Value = "'3.05"
$ws.Range("E30").Value = "  -1.35%  "

$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.144"
$ws.Range("E31").Value = "  +0.22%  "

$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'11.53"
$ws.Range("E32").Value = "  -0.31%  "

$ws.Range("E33").Value = "  -0.02%  "

$ws.Range("E34").Value = "  +0.37%  "

$ws.Range("E35").Value = "  -2.81%  "

$ws.Range("D36").Value = "'31.92"
$ws.Range("E36").Value = "  +4.93%  "

$ws.Range("D37").Value = "'1.67"
$ws.Range("E37").Value = "  +12.06%  "

$ws.Range("D38").Value = "'0.558"
$ws.Range("E38").Value = "  -1.99%  "

$ws.Range("D39").Value = "'8.48"
$ws.Range("E39").Value = "  +7.44%  "

$ws.Range("D40").Value = "'581.51"
$ws.Range("E40").Value = "  +1.29%  "

$ws.Range("E41").Value = "  +0.03%  "

$ws.Range("E42").Value = "  -0.74%  "

$ws.Range("D43").Value = "'0.908"
$ws.Range("E43").Value = "  -1.77%  "

$ws.Range("D44").Value = "'1.83"
$ws.Range("E44").Value = "  +4.26%  "

$ws.Range("D45").Value = "'5.74"
$ws.Range("E45").Value = "  +1.53%  "

$ws.Range("D46").Value = "'34.47"
$ws.Range("E46").Value = "  +31.89%  "

$ws.Range("D47").Value = "'2.28"
$ws.Range("E47").Value = "  +3.67%  "

$ws.Range("D49").Value = "'0.0415"
$ws.Range("E49").Value = "  -3.70%  "

$ws.Range("E50").Value = "  +1.06%  "

$ws.Range("D51").Value = "'8.18"
$ws.Range("E51").Value = "  +0.72%  "
